$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix cognacy of "north wind" units: correct L1 (col A) to match the actual
# language of Word1 (col B), move the old (incorrect) L1 value into L2 (col D),
# round the Levenshtein Distance (col G) to 6 decimal places, and clear the
# Cognates flag (col H) for word pairs that are not actually cognate.

$ws.Cells.Item(2322, 1).Value = "RU"
$ws.Cells.Item(2322, 4).Value = "BE"
$ws.Cells.Item(2322, 7).Value = 0.647059

$ws.Cells.Item(2323, 1).Value = "RU"
$ws.Cells.Item(2323, 4).Value = "UK"
$ws.Cells.Item(2323, 7).Value = 0.705882

$ws.Cells.Item(2324, 1).Value = "RU"
$ws.Cells.Item(2324, 4).Value = "PL"
$ws.Cells.Item(2324, 7).Value = 0.705882

$ws.Cells.Item(2325, 1).Value = "RU"
$ws.Cells.Item(2325, 4).Value = "CS"
$ws.Cells.Item(2325, 7).Value = 0.823529

$ws.Cells.Item(2326, 1).Value = "RU"
$ws.Cells.Item(2326, 4).Value = "SK"
$ws.Cells.Item(2326, 7).Value = 0.882353

$ws.Cells.Item(2327, 1).Value = "RU"
$ws.Cells.Item(2327, 4).Value = "SL"
$ws.Cells.Item(2327, 7).Value = 0.764706

$ws.Cells.Item(2328, 1).Value = "RU"
$ws.Cells.Item(2328, 4).Value = "HR"
$ws.Cells.Item(2328, 7).Value = 0.588235

$ws.Cells.Item(2329, 1).Value = "RU"
$ws.Cells.Item(2329, 4).Value = "BG"
$ws.Cells.Item(2329, 7).Value = 0.647059

$ws.Cells.Item(2330, 1).Value = "BE"
$ws.Cells.Item(2330, 4).Value = "UK"
$ws.Cells.Item(2330, 7).Value = 0.705882

$ws.Cells.Item(2331, 1).Value = "BE"
$ws.Cells.Item(2331, 4).Value = "PL"
$ws.Cells.Item(2331, 7).Value = 0.533333

$ws.Cells.Item(2332, 1).Value = "BE"
$ws.Cells.Item(2332, 4).Value = "CS"
$ws.Cells.Item(2332, 7).Value = 0.933333
$ws.Cells.Item(2332, 8).ClearContents()

$ws.Cells.Item(2333, 1).Value = "BE"
$ws.Cells.Item(2333, 4).Value = "SK"
$ws.Cells.Item(2333, 7).Value = 1.0
$ws.Cells.Item(2333, 8).ClearContents()

$ws.Cells.Item(2334, 1).Value = "BE"
$ws.Cells.Item(2334, 4).Value = "SL"
$ws.Cells.Item(2334, 7).Value = 0.8

$ws.Cells.Item(2335, 1).Value = "BE"
$ws.Cells.Item(2335, 4).Value = "HR"
$ws.Cells.Item(2335, 7).Value = 0.8125

$ws.Cells.Item(2336, 1).Value = "BE"
$ws.Cells.Item(2336, 4).Value = "BG"
$ws.Cells.Item(2336, 7).Value = 0.875

$ws.Cells.Item(2337, 1).Value = "UK"
$ws.Cells.Item(2337, 4).Value = "PL"
$ws.Cells.Item(2337, 7).Value = 0.705882

$ws.Cells.Item(2338, 1).Value = "UK"
$ws.Cells.Item(2338, 4).Value = "CS"
$ws.Cells.Item(2338, 7).Value = 1.0
$ws.Cells.Item(2338, 8).ClearContents()

$ws.Cells.Item(2339, 1).Value = "UK"
$ws.Cells.Item(2339, 4).Value = "SK"
$ws.Cells.Item(2339, 7).Value = 0.941176
$ws.Cells.Item(2339, 8).ClearContents()

$ws.Cells.Item(2340, 1).Value = "UK"
$ws.Cells.Item(2340, 4).Value = "SL"
$ws.Cells.Item(2340, 7).Value = 0.764706

$ws.Cells.Item(2341, 1).Value = "UK"
$ws.Cells.Item(2341, 4).Value = "HR"
$ws.Cells.Item(2341, 7).Value = 0.705882

$ws.Cells.Item(2342, 1).Value = "UK"
$ws.Cells.Item(2342, 4).Value = "BG"
$ws.Cells.Item(2342, 7).Value = 0.823529

$ws.Cells.Item(2343, 1).Value = "PL"
$ws.Cells.Item(2343, 4).Value = "CS"
$ws.Cells.Item(2343, 7).Value = 0.923077
$ws.Cells.Item(2343, 8).ClearContents()

$ws.Cells.Item(2344, 1).Value = "PL"
$ws.Cells.Item(2344, 4).Value = "SK"
$ws.Cells.Item(2344, 7).Value = 0.923077
$ws.Cells.Item(2344, 8).ClearContents()

$ws.Cells.Item(2345, 1).Value = "PL"
$ws.Cells.Item(2345, 4).Value = "SL"
$ws.Cells.Item(2345, 7).Value = 0.857143

$ws.Cells.Item(2346, 1).Value = "PL"
$ws.Cells.Item(2346, 4).Value = "HR"
$ws.Cells.Item(2346, 7).Value = 0.75

$ws.Cells.Item(2347, 1).Value = "PL"
$ws.Cells.Item(2347, 4).Value = "BG"
$ws.Cells.Item(2347, 7).Value = 0.625

$ws.Cells.Item(2348, 1).Value = "CS"
$ws.Cells.Item(2348, 4).Value = "SK"
$ws.Cells.Item(2348, 7).Value = 0.125

$ws.Cells.Item(2349, 1).Value = "CS"
$ws.Cells.Item(2349, 4).Value = "SL"
$ws.Cells.Item(2349, 7).Value = 0.857143

$ws.Cells.Item(2350, 1).Value = "CS"
$ws.Cells.Item(2350, 4).Value = "HR"
$ws.Cells.Item(2350, 7).Value = 0.8125

$ws.Cells.Item(2351, 1).Value = "CS"
$ws.Cells.Item(2351, 4).Value = "BG"
$ws.Cells.Item(2351, 7).Value = 0.625

$ws.Cells.Item(2352, 1).Value = "SK"
$ws.Cells.Item(2352, 4).Value = "SL"
$ws.Cells.Item(2352, 7).Value = 0.785714

$ws.Cells.Item(2353, 1).Value = "SK"
$ws.Cells.Item(2353, 4).Value = "HR"
$ws.Cells.Item(2353, 7).Value = 0.75

$ws.Cells.Item(2354, 1).Value = "SK"
$ws.Cells.Item(2354, 4).Value = "BG"
$ws.Cells.Item(2354, 7).Value = 0.6875

$ws.Cells.Item(2355, 1).Value = "SL"
$ws.Cells.Item(2355, 4).Value = "HR"
$ws.Cells.Item(2355, 7).Value = 0.625

$ws.Cells.Item(2356, 1).Value = "SL"
$ws.Cells.Item(2356, 4).Value = "BG"
$ws.Cells.Item(2356, 7).Value = 0.8125

$ws.Cells.Item(2357, 1).Value = "HR"
$ws.Cells.Item(2357, 4).Value = "BG"
$ws.Cells.Item(2357, 7).Value = 0.6875

$ws.Cells.Item(2358, 1).Value = "RU"
$ws.Cells.Item(2358, 4).Value = "BE"
$ws.Cells.Item(2358, 7).Value = 0.647059

$ws.Cells.Item(2359, 1).Value = "RU"
$ws.Cells.Item(2359, 4).Value = "UK"
$ws.Cells.Item(2359, 7).Value = 0.705882

$ws.Cells.Item(2360, 1).Value = "RU"
$ws.Cells.Item(2360, 4).Value = "PL"
$ws.Cells.Item(2360, 7).Value = 0.705882

$ws.Cells.Item(2361, 1).Value = "RU"
$ws.Cells.Item(2361, 4).Value = "CS"
$ws.Cells.Item(2361, 7).Value = 0.823529

$ws.Cells.Item(2362, 1).Value = "RU"
$ws.Cells.Item(2362, 4).Value = "SK"
$ws.Cells.Item(2362, 7).Value = 0.882353

$ws.Cells.Item(2363, 1).Value = "RU"
$ws.Cells.Item(2363, 4).Value = "SL"
$ws.Cells.Item(2363, 7).Value = 0.764706

$ws.Cells.Item(2364, 1).Value = "RU"
$ws.Cells.Item(2364, 4).Value = "HR"
$ws.Cells.Item(2364, 7).Value = 0.588235

$ws.Cells.Item(2365, 1).Value = "RU"
$ws.Cells.Item(2365, 4).Value = "BG"
$ws.Cells.Item(2365, 7).Value = 0.647059

$ws.Cells.Item(2366, 1).Value = "BE"
$ws.Cells.Item(2366, 4).Value = "UK"
$ws.Cells.Item(2366, 7).Value = 0.764706

$ws.Cells.Item(2367, 1).Value = "BE"
$ws.Cells.Item(2367, 4).Value = "PL"
$ws.Cells.Item(2367, 7).Value = 0.533333

$ws.Cells.Item(2368, 1).Value = "BE"
$ws.Cells.Item(2368, 4).Value = "CS"
$ws.Cells.Item(2368, 7).Value = 0.933333
$ws.Cells.Item(2368, 8).ClearContents()

$ws.Cells.Item(2369, 1).Value = "BE"
$ws.Cells.Item(2369, 4).Value = "SK"
$ws.Cells.Item(2369, 7).Value = 1.0
$ws.Cells.Item(2369, 8).ClearContents()

$ws.Cells.Item(2370, 1).Value = "BE"
$ws.Cells.Item(2370, 4).Value = "SL"
$ws.Cells.Item(2370, 7).Value = 0.8

$ws.Cells.Item(2371, 1).Value = "BE"
$ws.Cells.Item(2371, 4).Value = "HR"
$ws.Cells.Item(2371, 7).Value = 0.8125

$ws.Cells.Item(2372, 1).Value = "BE"
$ws.Cells.Item(2372, 4).Value = "BG"
$ws.Cells.Item(2372, 7).Value = 0.875

$ws.Cells.Item(2373, 1).Value = "UK"
$ws.Cells.Item(2373, 4).Value = "PL"
$ws.Cells.Item(2373, 7).Value = 0.705882

$ws.Cells.Item(2374, 1).Value = "UK"
$ws.Cells.Item(2374, 4).Value = "CS"
$ws.Cells.Item(2374, 7).Value = 1.0
$ws.Cells.Item(2374, 8).ClearContents()

$ws.Cells.Item(2375, 1).Value = "UK"
$ws.Cells.Item(2375, 4).Value = "SK"
$ws.Cells.Item(2375, 7).Value = 0.941176
$ws.Cells.Item(2375, 8).ClearContents()

$ws.Cells.Item(2376, 1).Value = "UK"
$ws.Cells.Item(2376, 4).Value = "SL"
$ws.Cells.Item(2376, 7).Value = 0.764706

$ws.Cells.Item(2377, 1).Value = "UK"
$ws.Cells.Item(2377, 4).Value = "HR"
$ws.Cells.Item(2377, 7).Value = 0.705882

$ws.Cells.Item(2378, 1).Value = "UK"
$ws.Cells.Item(2378, 4).Value = "BG"
$ws.Cells.Item(2378, 7).Value = 0.823529

$ws.Cells.Item(2379, 1).Value = "PL"
$ws.Cells.Item(2379, 4).Value = "CS"
$ws.Cells.Item(2379, 7).Value = 0.923077
$ws.Cells.Item(2379, 8).ClearContents()

$ws.Cells.Item(2380, 1).Value = "PL"
$ws.Cells.Item(2380, 4).Value = "SK"
$ws.Cells.Item(2380, 7).Value = 0.923077
$ws.Cells.Item(2380, 8).ClearContents()

$ws.Cells.Item(2381, 1).Value = "PL"
$ws.Cells.Item(2381, 4).Value = "SL"
$ws.Cells.Item(2381, 7).Value = 0.857143

$ws.Cells.Item(2382, 1).Value = "PL"
$ws.Cells.Item(2382, 4).Value = "HR"
$ws.Cells.Item(2382, 7).Value = 0.75

$ws.Cells.Item(2383, 1).Value = "PL"
$ws.Cells.Item(2383, 4).Value = "BG"
$ws.Cells.Item(2383, 7).Value = 0.625

$ws.Cells.Item(2384, 1).Value = "CS"
$ws.Cells.Item(2384, 4).Value = "SK"
$ws.Cells.Item(2384, 7).Value = 0.125

$ws.Cells.Item(2385, 1).Value = "CS"
$ws.Cells.Item(2385, 4).Value = "SL"
$ws.Cells.Item(2385, 7).Value = 0.857143

$ws.Cells.Item(2386, 1).Value = "CS"
$ws.Cells.Item(2386, 4).Value = "HR"
$ws.Cells.Item(2386, 7).Value = 0.8125

$ws.Cells.Item(2387, 1).Value = "CS"
$ws.Cells.Item(2387, 4).Value = "BG"
$ws.Cells.Item(2387, 7).Value = 0.625

$ws.Cells.Item(2388, 1).Value = "SK"
$ws.Cells.Item(2388, 4).Value = "SL"
$ws.Cells.Item(2388, 7).Value = 0.785714

$ws.Cells.Item(2389, 1).Value = "SK"
$ws.Cells.Item(2389, 4).Value = "HR"
$ws.Cells.Item(2389, 7).Value = 0.75

$ws.Cells.Item(2390, 1).Value = "SK"
$ws.Cells.Item(2390, 4).Value = "BG"
$ws.Cells.Item(2390, 7).Value = 0.6875

$ws.Cells.Item(2391, 1).Value = "SL"
$ws.Cells.Item(2391, 4).Value = "HR"
$ws.Cells.Item(2391, 7).Value = 0.625

$ws.Cells.Item(2392, 1).Value = "SL"
$ws.Cells.Item(2392, 4).Value = "BG"
$ws.Cells.Item(2392, 7).Value = 0.8125

$ws.Cells.Item(2393, 1).Value = "HR"
$ws.Cells.Item(2393, 4).Value = "BG"
$ws.Cells.Item(2393, 7).Value = 0.6875

$ws.Cells.Item(2394, 1).Value = "RU"
$ws.Cells.Item(2394, 4).Value = "BE"
$ws.Cells.Item(2394, 7).Value = 0.647059

$ws.Cells.Item(2395, 1).Value = "RU"
$ws.Cells.Item(2395, 4).Value = "UK"
$ws.Cells.Item(2395, 7).Value = 0.705882

$ws.Cells.Item(2396, 1).Value = "RU"
$ws.Cells.Item(2396, 4).Value = "PL"
$ws.Cells.Item(2396, 7).Value = 0.705882

$ws.Cells.Item(2397, 1).Value = "RU"
$ws.Cells.Item(2397, 4).Value = "CS"
$ws.Cells.Item(2397, 7).Value = 0.823529

$ws.Cells.Item(2398, 1).Value = "RU"
$ws.Cells.Item(2398, 4).Value = "SK"
$ws.Cells.Item(2398, 7).Value = 0.882353

$ws.Cells.Item(2399, 1).Value = "RU"
$ws.Cells.Item(2399, 4).Value = "SL"
$ws.Cells.Item(2399, 7).Value = 0.764706

$ws.Cells.Item(2400, 1).Value = "RU"
$ws.Cells.Item(2400, 4).Value = "HR"
$ws.Cells.Item(2400, 7).Value = 0.588235

$ws.Cells.Item(2401, 1).Value = "RU"
$ws.Cells.Item(2401, 4).Value = "BG"
$ws.Cells.Item(2401, 7).Value = 0.647059

$ws.Cells.Item(2402, 1).Value = "BE"
$ws.Cells.Item(2402, 4).Value = "UK"
$ws.Cells.Item(2402, 7).Value = 0.705882

$ws.Cells.Item(2403, 1).Value = "BE"
$ws.Cells.Item(2403, 4).Value = "PL"
$ws.Cells.Item(2403, 7).Value = 0.533333

$ws.Cells.Item(2404, 1).Value = "BE"
$ws.Cells.Item(2404, 4).Value = "CS"
$ws.Cells.Item(2404, 7).Value = 0.933333
$ws.Cells.Item(2404, 8).ClearContents()

$ws.Cells.Item(2405, 1).Value = "BE"
$ws.Cells.Item(2405, 4).Value = "SK"
$ws.Cells.Item(2405, 7).Value = 1.0
$ws.Cells.Item(2405, 8).ClearContents()

$ws.Cells.Item(2406, 1).Value = "BE"
$ws.Cells.Item(2406, 4).Value = "SL"
$ws.Cells.Item(2406, 7).Value = 0.8

$ws.Cells.Item(2407, 1).Value = "BE"
$ws.Cells.Item(2407, 4).Value = "HR"
$ws.Cells.Item(2407, 7).Value = 0.8125

$ws.Cells.Item(2408, 1).Value = "BE"
$ws.Cells.Item(2408, 4).Value = "BG"
$ws.Cells.Item(2408, 7).Value = 0.875

$ws.Cells.Item(2409, 1).Value = "UK"
$ws.Cells.Item(2409, 4).Value = "PL"
$ws.Cells.Item(2409, 7).Value = 0.705882

$ws.Cells.Item(2410, 1).Value = "UK"
$ws.Cells.Item(2410, 4).Value = "CS"
$ws.Cells.Item(2410, 7).Value = 1.0
$ws.Cells.Item(2410, 8).ClearContents()

$ws.Cells.Item(2411, 1).Value = "UK"
$ws.Cells.Item(2411, 4).Value = "SK"
$ws.Cells.Item(2411, 7).Value = 0.941176
$ws.Cells.Item(2411, 8).ClearContents()

$ws.Cells.Item(2412, 1).Value = "UK"
$ws.Cells.Item(2412, 4).Value = "SL"
$ws.Cells.Item(2412, 7).Value = 0.764706

$ws.Cells.Item(2413, 1).Value = "UK"
$ws.Cells.Item(2413, 4).Value = "HR"
$ws.Cells.Item(2413, 7).Value = 0.705882

$ws.Cells.Item(2414, 1).Value = "UK"
$ws.Cells.Item(2414, 4).Value = "BG"
$ws.Cells.Item(2414, 7).Value = 0.823529

$ws.Cells.Item(2415, 1).Value = "PL"
$ws.Cells.Item(2415, 4).Value = "CS"
$ws.Cells.Item(2415, 7).Value = 0.923077
$ws.Cells.Item(2415, 8).ClearContents()

$ws.Cells.Item(2416, 1).Value = "PL"
$ws.Cells.Item(2416, 4).Value = "SK"
$ws.Cells.Item(2416, 7).Value = 0.923077
$ws.Cells.Item(2416, 8).ClearContents()

$ws.Cells.Item(2417, 1).Value = "PL"
$ws.Cells.Item(2417, 4).Value = "SL"
$ws.Cells.Item(2417, 7).Value = 0.857143

$ws.Cells.Item(2418, 1).Value = "PL"
$ws.Cells.Item(2418, 4).Value = "HR"
$ws.Cells.Item(2418, 7).Value = 0.75

$ws.Cells.Item(2419, 1).Value = "PL"
$ws.Cells.Item(2419, 4).Value = "BG"
$ws.Cells.Item(2419, 7).Value = 0.625

$ws.Cells.Item(2420, 1).Value = "CS"
$ws.Cells.Item(2420, 4).Value = "SK"
$ws.Cells.Item(2420, 7).Value = 0.25

$ws.Cells.Item(2421, 1).Value = "CS"
$ws.Cells.Item(2421, 4).Value = "SL"
$ws.Cells.Item(2421, 7).Value = 0.857143

$ws.Cells.Item(2422, 1).Value = "CS"
$ws.Cells.Item(2422, 4).Value = "HR"
$ws.Cells.Item(2422, 7).Value = 0.8125

$ws.Cells.Item(2423, 1).Value = "CS"
$ws.Cells.Item(2423, 4).Value = "BG"
$ws.Cells.Item(2423, 7).Value = 0.625

$ws.Cells.Item(2424, 1).Value = "SK"
$ws.Cells.Item(2424, 4).Value = "SL"
$ws.Cells.Item(2424, 7).Value = 0.785714

$ws.Cells.Item(2425, 1).Value = "SK"
$ws.Cells.Item(2425, 4).Value = "HR"
$ws.Cells.Item(2425, 7).Value = 0.75

$ws.Cells.Item(2426, 1).Value = "SK"
$ws.Cells.Item(2426, 4).Value = "BG"
$ws.Cells.Item(2426, 7).Value = 0.6875

$ws.Cells.Item(2427, 1).Value = "SL"
$ws.Cells.Item(2427, 4).Value = "HR"
$ws.Cells.Item(2427, 7).Value = 0.625

$ws.Cells.Item(2428, 1).Value = "SL"
$ws.Cells.Item(2428, 4).Value = "BG"
$ws.Cells.Item(2428, 7).Value = 0.8125

$ws.Cells.Item(2429, 1).Value = "HR"
$ws.Cells.Item(2429, 4).Value = "BG"
$ws.Cells.Item(2429, 7).Value = 0.6875

$ws.Cells.Item(2430, 1).Value = "RU"
$ws.Cells.Item(2430, 4).Value = "BE"
$ws.Cells.Item(2430, 7).Value = 0.647059

$ws.Cells.Item(2431, 1).Value = "RU"
$ws.Cells.Item(2431, 4).Value = "UK"
$ws.Cells.Item(2431, 7).Value = 0.705882

$ws.Cells.Item(2432, 1).Value = "RU"
$ws.Cells.Item(2432, 4).Value = "PL"
$ws.Cells.Item(2432, 7).Value = 0.705882

$ws.Cells.Item(2433, 1).Value = "RU"
$ws.Cells.Item(2433, 4).Value = "CS"
$ws.Cells.Item(2433, 7).Value = 0.823529

$ws.Cells.Item(2434, 1).Value = "RU"
$ws.Cells.Item(2434, 4).Value = "SK"
$ws.Cells.Item(2434, 7).Value = 0.882353

$ws.Cells.Item(2435, 1).Value = "RU"
$ws.Cells.Item(2435, 4).Value = "SL"
$ws.Cells.Item(2435, 7).Value = 0.764706

$ws.Cells.Item(2436, 1).Value = "RU"
$ws.Cells.Item(2436, 4).Value = "HR"
$ws.Cells.Item(2436, 7).Value = 0.588235

$ws.Cells.Item(2437, 1).Value = "RU"
$ws.Cells.Item(2437, 4).Value = "BG"
$ws.Cells.Item(2437, 7).Value = 0.647059

$ws.Cells.Item(2438, 1).Value = "BE"
$ws.Cells.Item(2438, 4).Value = "UK"
$ws.Cells.Item(2438, 7).Value = 0.705882

$ws.Cells.Item(2439, 1).Value = "BE"
$ws.Cells.Item(2439, 4).Value = "PL"
$ws.Cells.Item(2439, 7).Value = 0.533333

$ws.Cells.Item(2440, 1).Value = "BE"
$ws.Cells.Item(2440, 4).Value = "CS"
$ws.Cells.Item(2440, 7).Value = 0.933333
$ws.Cells.Item(2440, 8).ClearContents()

$ws.Cells.Item(2441, 1).Value = "BE"
$ws.Cells.Item(2441, 4).Value = "SK"
$ws.Cells.Item(2441, 7).Value = 1.0
$ws.Cells.Item(2441, 8).ClearContents()

$ws.Cells.Item(2442, 1).Value = "BE"
$ws.Cells.Item(2442, 4).Value = "SL"
$ws.Cells.Item(2442, 7).Value = 0.8

$ws.Cells.Item(2443, 1).Value = "BE"
$ws.Cells.Item(2443, 4).Value = "HR"
$ws.Cells.Item(2443, 7).Value = 0.8125

$ws.Cells.Item(2444, 1).Value = "BE"
$ws.Cells.Item(2444, 4).Value = "BG"
$ws.Cells.Item(2444, 7).Value = 0.875

$ws.Cells.Item(2445, 1).Value = "UK"
$ws.Cells.Item(2445, 4).Value = "PL"
$ws.Cells.Item(2445, 7).Value = 0.705882

$ws.Cells.Item(2446, 1).Value = "UK"
$ws.Cells.Item(2446, 4).Value = "CS"
$ws.Cells.Item(2446, 7).Value = 1.0
$ws.Cells.Item(2446, 8).ClearContents()

$ws.Cells.Item(2447, 1).Value = "UK"
$ws.Cells.Item(2447, 4).Value = "SK"
$ws.Cells.Item(2447, 7).Value = 0.941176
$ws.Cells.Item(2447, 8).ClearContents()

$ws.Cells.Item(2448, 1).Value = "UK"
$ws.Cells.Item(2448, 4).Value = "SL"
$ws.Cells.Item(2448, 7).Value = 0.764706

$ws.Cells.Item(2449, 1).Value = "UK"
$ws.Cells.Item(2449, 4).Value = "HR"
$ws.Cells.Item(2449, 7).Value = 0.705882

$ws.Cells.Item(2450, 1).Value = "UK"
$ws.Cells.Item(2450, 4).Value = "BG"
$ws.Cells.Item(2450, 7).Value = 0.823529

$ws.Cells.Item(2451, 1).Value = "PL"
$ws.Cells.Item(2451, 4).Value = "CS"
$ws.Cells.Item(2451, 7).Value = 0.923077
$ws.Cells.Item(2451, 8).ClearContents()

$ws.Cells.Item(2452, 1).Value = "PL"
$ws.Cells.Item(2452, 4).Value = "SK"
$ws.Cells.Item(2452, 7).Value = 0.923077
$ws.Cells.Item(2452, 8).ClearContents()

$ws.Cells.Item(2453, 1).Value = "PL"
$ws.Cells.Item(2453, 4).Value = "SL"
$ws.Cells.Item(2453, 7).Value = 0.857143

$ws.Cells.Item(2454, 1).Value = "PL"
$ws.Cells.Item(2454, 4).Value = "HR"
$ws.Cells.Item(2454, 7).Value = 0.75

$ws.Cells.Item(2455, 1).Value = "PL"
$ws.Cells.Item(2455, 4).Value = "BG"
$ws.Cells.Item(2455, 7).Value = 0.625

$ws.Cells.Item(2456, 1).Value = "CS"
$ws.Cells.Item(2456, 4).Value = "SK"
$ws.Cells.Item(2456, 7).Value = 0.125

$ws.Cells.Item(2457, 1).Value = "CS"
$ws.Cells.Item(2457, 4).Value = "SL"
$ws.Cells.Item(2457, 7).Value = 0.857143

$ws.Cells.Item(2458, 1).Value = "CS"
$ws.Cells.Item(2458, 4).Value = "HR"
$ws.Cells.Item(2458, 7).Value = 0.8125

$ws.Cells.Item(2459, 1).Value = "CS"
$ws.Cells.Item(2459, 4).Value = "BG"
$ws.Cells.Item(2459, 7).Value = 0.625

$ws.Cells.Item(2460, 1).Value = "SK"
$ws.Cells.Item(2460, 4).Value = "SL"
$ws.Cells.Item(2460, 7).Value = 0.785714

$ws.Cells.Item(2461, 1).Value = "SK"
$ws.Cells.Item(2461, 4).Value = "HR"
$ws.Cells.Item(2461, 7).Value = 0.75

$ws.Cells.Item(2462, 1).Value = "SK"
$ws.Cells.Item(2462, 4).Value = "BG"
$ws.Cells.Item(2462, 7).Value = 0.6875

$ws.Cells.Item(2463, 1).Value = "SL"
$ws.Cells.Item(2463, 4).Value = "HR"
$ws.Cells.Item(2463, 7).Value = 0.625

$ws.Cells.Item(2464, 1).Value = "SL"
$ws.Cells.Item(2464, 4).Value = "BG"
$ws.Cells.Item(2464, 7).Value = 0.8125

$ws.Cells.Item(2465, 1).Value = "HR"
$ws.Cells.Item(2465, 4).Value = "BG"
$ws.Cells.Item(2465, 7).Value = 0.6875
